$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 0.4980883136293128;  C = 0.9900820039247082;  D = 0.5772119266417263;  G = 1.558460351833249; H = 0.9990000000000001 }
    3  = @{ B = 0.239056291290989;   C = 0.9953276387749029;  D = 0.3843433050552175;  G = 1.558460351833249; H = 0.9990000000000001 }
    4  = @{ B = 0.3424301550227666;  C = 0.9934124052592661;  D = 0.4655874903828661;  G = 1.558460351833249; H = 0.9990000000000001 }
    5  = @{ B = 0.4178647680449223;  C = 0.9917602788820304;  D = 0.5037273826483598;  G = 1.558460351833249; H = 0.9990000000000001 }
    6  = @{ B = 0.7195241642125411;  C = 0.9788656904573576;  D = 0.6549930491151482;  G = 1.558460351833249; H = 0.9990000000000001 }
    7  = @{ B = 0.3829998621080163;  C = 0.9947068382029508;  D = 0.5150875870527777;  G = 1.558460351833249; H = 0.9990000000000001 }
    8  = @{ B = 0.1529799275151037;  C = 0.9984232685064792;  D = 0.3311603784071903;  G = 1.558460351833249; H = 0.9990000000000001 }
    9  = @{ B = 0.4978591271369693;  C = 0.9970343952301695;  D = 0.5844367220048292;  G = 1.558460351833249; H = 0.9990000000000001 }
    10 = @{ B = 0.1044224412242731;  C = 0.9980924492926455;  D = 0.2348577086854635;  G = 1.558460351833249; H = 0.9990000000000001 }
    11 = @{ B = 0.3134948840143841;  C = 0.97683785591566;    D = 0.4293759673639526;  G = 1.558460351833249; H = 0.9990000000000001 }
    12 = @{ B = 0.05243601867023975; C = 0.9984579461493137;  D = 0.1633510054277421;  G = 1.558460351833249; H = 0.9990000000000001 }
    13 = @{ B = 0.07632337034479625; C = 0.9992760386281136;  D = 0.2081515299792747;  G = 1.558460351833249; H = 0.9990000000000001 }
    14 = @{ B = 0.08726761627618847; C = 0.998818450344211;   D = 0.2390931973700207;  G = 1.558460351833249; H = 0.9990000000000001 }
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("B$row").Value = $rowValues.B
    $ws.Range("C$row").Value = $rowValues.C
    $ws.Range("D$row").Value = $rowValues.D
    $ws.Range("G$row").Value = $rowValues.G
    $ws.Range("H$row").Value = $rowValues.H
}
